$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")
$ws.Rows.Item(2).Delete()
